$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.449.25'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '1.726.61'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.33%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'244.62"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("E6").Value = '  +0.29%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.4788"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  +1.28%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.2684"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  +1.23%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '1.730.52'
$ws.Range("E10").Value = '  +1.02%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.07129"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("E12").Value = '  +2.53%  '
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.6168"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = '  +4.06%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'4.539"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  +2.47%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").Value = "'77.19"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '26.460.82'
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("E18").Value = '  +0.27%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'0.000006937"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  +1.82%  '
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'11.69"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '1.953.02'
$ws.Range("E21").Value = '  +2.34%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'4.541"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -1.27%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'8.914"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.53%  '
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'5.299"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = '  -1.12%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'136.47"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.57%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'15.35"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  +1.03%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").Value = "'1.801"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  +2.16%  '
$ws.Range("E28").Value = '  +0.19%  '
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'106.71"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  -0.01%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'3.979"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  -1.26%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'0.08011"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  +3.06%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.730"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.63%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'0.04556"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +2.94%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'0.9997"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  +0.41%  '
$ws.Range("E35").Value = '  +0.48%  '
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.6361"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = '  +1.84%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.9884"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +1.23%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.9337"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  +1.16%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'2.041"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  +6.23%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.411"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +0.03%  '
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").Value = "'107.08"
$ws.Range("D41").Style = $origStyle
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").Value = "'1.003"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  +1.42%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").Value = "'5.603"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +9.14%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.3904"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +2.18%  '
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'7.029"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  +12.05%  '
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'0.1191"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = '  +4.16%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.05322"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'7.919"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = "'30.96"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +0.50%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'1.268"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +3.29%  '
